$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying stat cells (runs, balls, fours, sixes) are stored as text,
# so force text formatting before writing the new values to keep them as
# text rather than being auto-converted to numbers.
$statRange = $ws.Range("C2:F4")
$statRange.NumberFormat = "@"

# Rotate the three player-innings rows (2-4) so that:
#   new row 2 = old row 3
#   new row 3 = old row 4
#   new row 4 = old row 2
$ws.Range("C2").Value = "6"
$ws.Range("D2").Value = "6"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "0"

$ws.Range("C3").Value = "34"
$ws.Range("D3").Value = "36"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "2"

$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "2"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"
